$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 64: 2025-11-01 data for 四方坪站
$ws.Range("A64").Value = 45962
$ws.Range("B64").Value = "四方坪站"
$ws.Range("C64").Value = 8560.7999999999993
$ws.Range("D64").Value = 7444.11
$ws.Range("E64").Value = 2843.51
$ws.Range("F64").Value = 382

# New row 65: 2025-11-01 data for 高岭站
$ws.Range("A65").Value = 45962
$ws.Range("B65").Value = "高岭站"
$ws.Range("C65").Value = 3801.21
$ws.Range("D65").Value = 3298.57
$ws.Range("E65").Value = 934.44
$ws.Range("F65").Value = 158

# Update the selected/active cell to match the new extent of data
$ws.Range("K65").Select()
